# This script updates the "Comments" list on Sheet1:
#  - Removes the out-of-place "this is test message" comment that was
#    sitting at row 8, which shifts the subsequent comments up by one row.
#  - Adds a new comment "Wow that's great" right after "hlo sir".
#  - Re-adds "this is test message" further down the list, just before "hi"
#    (immediately after "Interesting concept, curious to review").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "this is test message" row (currently row 8).
$ws.Rows.Item(8).Delete()

# Insert a new row for "Wow that's great" right after "hlo sir" (now row 8).
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "Wow that's great"

# Insert a new row to re-add "this is test message" right before "hi"
# (which is now at row 49 after the delete/insert above).
$ws.Rows.Item(49).Insert()
$ws.Range("A49").Value = "this is test message"
